$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9663489460945129
$ws.Range("B1").Value = 0.7534570693969727
$ws.Range("C1").Value = 4.007230758666992
$ws.Range("D1").Value = 2.977682590484619
$ws.Range("E1").Value = 0.8051601052284241
